# Auto-generated Excel COM-interop script to apply value updates
# per the diff between before.xlsx and the target workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 735.875
$ws.Range("I2").Value = 881.5
$ws.Range("J2").Value = 299
$ws.Range("K2").Value = 881.5
$ws.Range("L2").Value = 299
$ws.Range("M2").Value = -768.5
$ws.Range("N2").Value = -525
$ws.Range("H17").Value = 1082661
$ws.Range("J17").Value = 1082661
$ws.Range("L17").Value = 3247983
$ws.Range("N17").Value = -3248319
$ws.Range("H33").Value = 432.7857
$ws.Range("I33").Value = 158.23077
$ws.Range("K33").Value = 158.23077
$ws.Range("M33").Value = 70.76922999999999
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("N76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("N79").Value = 0
$ws.Range("H86").Value = 48700.555
$ws.Range("I86").Value = 951.25
$ws.Range("J86").Value = 86900
$ws.Range("K86").Value = 951.25
$ws.Range("L86").Value = 86900
$ws.Range("M86").Value = 171.75
$ws.Range("N86").Value = -89146
$ws.Range("H89").Value = 48700.555
$ws.Range("I89").Value = 951.25
$ws.Range("J89").Value = 86900
$ws.Range("K89").Value = 4756.25
$ws.Range("L89").Value = 434500
$ws.Range("M89").Value = 859.75
$ws.Range("N89").Value = -445732
$ws.Range("H132").Value = 3064.6326
$ws.Range("I132").Value = 3305.9744
$ws.Range("J132").Value = 2123.4
$ws.Range("K132").Value = 9917.923200000001
$ws.Range("L132").Value = 6370.200000000001
$ws.Range("M132").Value = -7387.923200000001
$ws.Range("N132").Value = -11430.2
$ws.Range("H138").Value = 4669.6885
$ws.Range("J138").Value = 6406.6343
$ws.Range("L138").Value = 19219.9029
$ws.Range("N138").Value = -29499.9029
$ws.Range("M76").ClearContents()
$ws.Range("M79").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1171828.9
$ws.Range("I122").Value = 3806
$ws.Range("K122").Value = 11418
$ws.Range("M122").Value = -8968
$ws.Range("H132").Value = 2750.697
$ws.Range("I132").Value = 2364.2173
$ws.Range("J132").Value = 3639.6
$ws.Range("K132").Value = 7092.651899999999
$ws.Range("L132").Value = 10918.8
$ws.Range("M132").Value = -4562.651899999999
$ws.Range("N132").Value = -15978.8
$ws.Range("H133").Value = 69997.25
$ws.Range("J133").Value = 69997.25
$ws.Range("L133").Value = 69997.25
$ws.Range("N133").Value = -75057.25
$ws.Range("H141").Value = 60429
$ws.Range("J141").Value = 60429
$ws.Range("L141").Value = 60429
$ws.Range("N141").Value = -70789

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3326.7827
$ws.Range("I107").Value = 2806.8333
$ws.Range("J107").Value = 5198.6
$ws.Range("K107").Value = 2806.8333
$ws.Range("L107").Value = 5198.6
$ws.Range("M107").Value = -886.8332999999998
$ws.Range("N107").Value = -9038.6
$ws.Range("H132").Value = 57629.668
$ws.Range("J132").Value = 57629.668
$ws.Range("L132").Value = 57629.668
$ws.Range("N132").Value = -67749.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 12324.966
$ws.Range("I107").Value = 13888.96
$ws.Range("K107").Value = 13888.96
$ws.Range("M107").Value = -11968.96
$ws.Range("H132").Value = 55049.7
$ws.Range("I132").Value = 18562.125
$ws.Range("K132").Value = 55686.375
$ws.Range("M132").Value = -53156.375
$ws.Range("H141").Value = 125080
$ws.Range("J141").Value = 125080
$ws.Range("L141").Value = 125080
$ws.Range("N141").Value = -135440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 277.35715
$ws.Range("I2").Value = 226.16667
$ws.Range("J2").Value = 315.75
$ws.Range("K2").Value = 1357.00002
$ws.Range("L2").Value = 1894.5
$ws.Range("M2").Value = -1244.00002
$ws.Range("N2").Value = -2120.5
$ws.Range("H23").Value = 11111253
$ws.Range("I23").Value = 56.22222
$ws.Range("J23").Value = 27778048
$ws.Range("K23").Value = 168.66666
$ws.Range("L23").Value = 83334144
$ws.Range("M23").Value = 66.33333999999999
$ws.Range("N23").Value = -83334614
$ws.Range("H34").Value = 2086905.9
$ws.Range("I34").Value = 2780041
$ws.Range("J34").Value = 7500
$ws.Range("K34").Value = 8340123
$ws.Range("L34").Value = 22500
$ws.Range("M34").Value = -8340039
$ws.Range("N34").Value = -22668
$ws.Range("H38").Value = 1336
$ws.Range("I38").Value = 237.38461
$ws.Range("J38").Value = 2288.1333
$ws.Range("K38").Value = 712.15383
$ws.Range("L38").Value = 6864.3999
$ws.Range("M38").Value = -365.15383
$ws.Range("N38").Value = -7558.3999
$ws.Range("H39").Value = 3961.75
$ws.Range("I39").Value = 283.33334
$ws.Range("J39").Value = 14997
$ws.Range("K39").Value = 850.0000200000001
$ws.Range("L39").Value = 44991
$ws.Range("M39").Value = -556.0000200000001
$ws.Range("N39").Value = -45579
$ws.Range("H55").Value = 6850.3076
$ws.Range("J55").Value = 10558.125
$ws.Range("L55").Value = 31674.375
$ws.Range("N55").Value = -32028.375
$ws.Range("H107").Value = 2271.9
$ws.Range("J107").Value = 2239.875
$ws.Range("L107").Value = 6719.625
$ws.Range("N107").Value = -10559.625
$ws.Range("H122").Value = 1000.9286
$ws.Range("J122").Value = 1143
$ws.Range("L122").Value = 10287
$ws.Range("N122").Value = -15187
$ws.Range("H131").Value = 14087101
$ws.Range("I131").Value = 111120240
$ws.Range("J131").Value = 1645.0483
$ws.Range("K131").Value = 333360720
$ws.Range("L131").Value = 4935.144899999999
$ws.Range("M131").Value = -333355680
$ws.Range("N131").Value = -15015.1449

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7571.0835
$ws.Range("I70").Value = 7182.0625
$ws.Range("J70").Value = 8349.125
$ws.Range("K70").Value = 7182.0625
$ws.Range("L70").Value = 8349.125
$ws.Range("M70").Value = -6912.0625
$ws.Range("N70").Value = -8889.125
$ws.Range("H73").Value = 7571.0835
$ws.Range("I73").Value = 7182.0625
$ws.Range("J73").Value = 8349.125
$ws.Range("K73").Value = 7182.0625
$ws.Range("L73").Value = 8349.125
$ws.Range("M73").Value = -6246.0625
$ws.Range("N73").Value = -10221.125
$ws.Range("H132").Value = 4158.533
$ws.Range("I132").Value = 2959.75
$ws.Range("J132").Value = 5528.5713
$ws.Range("K132").Value = 8879.25
$ws.Range("L132").Value = 16585.7139
$ws.Range("M132").Value = -6349.25
$ws.Range("N132").Value = -21645.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 280
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 78991.5
$ws.Range("J46").Value = 78991.5
$ws.Range("L46").Value = 78991.5
$ws.Range("N46").Value = -79453.5
$ws.Range("H134").Value = 78991.5
$ws.Range("J134").Value = 78991.5
$ws.Range("L134").Value = 236974.5
$ws.Range("N134").Value = -242044.5
$ws.Range("H136").Value = 4487.7334
$ws.Range("I136").Value = 3147.4443
$ws.Range("J136").Value = 6498.1665
$ws.Range("K136").Value = 9442.332900000001
$ws.Range("L136").Value = 19494.4995
$ws.Range("M136").Value = -6892.332900000001
$ws.Range("N136").Value = -24594.4995
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("N138").Value = 0
